$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84 - this shifts existing rows 84..181 down to 85..182
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record
$ws.Range("A84").Value = 10
$ws.Range("B84").Value = "Vega Modelo de Temuco"
$ws.Range("C84").Value = "La Araucanía"
$ws.Range("D84").Value = 44781
$ws.Range("E84").Value = 9
$ws.Range("F84").Value = 100114007
$ws.Range("G84").Value = "Jengibre"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 170
$ws.Range("K84").Value = 16000
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = 16941
$ws.Range("N84").Value = '$/caja 13 kilos'
$ws.Range("O84").Value = "Perú"
$ws.Range("P84").Value = 1303
$ws.Range("Q84").Value = 13
$ws.Range("R84").Value = "Hortaliza"
